# Apply the edits described by the commit:
# "New executable which accounts for user input"
#
# 1) B6 (areaPV) now computes capPV / etaPV_rated instead of costPV / capPV.
# 2) The active/selected cell on Sheet1 moves from B9 to G9.
# 3) Best-effort: nudge the workbook window position to match the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the formula in B6: areaPV = capPV / etaPV_rated  (B5 / B3)
$ws.Range("B6").Formula = "=B5/B3"

# Move the selection to G9 and make sure Sheet1 is the active sheet/window.
$ws.Activate()
$ws.Range("G9").Select()

# Best-effort: reflect the new window position recorded in the saved file.
$win = $excel.Windows.Item(1)
$win.Left = 5190
$win.Top = 3810
